# Append the latest file-usage snapshot row to the "Data" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Next empty row right after the current last row of data (row 84 -> 85)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$values = @(45816.677997685183, 14, 6, 414, 650, 622, 692, 5568, 692, 2, 2, 684, 30, 5916, 7100)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($newRow, $i + 1).Value = $values[$i]
}
